# Update line-loading results for the 380 kV case (rows 2-25, cols B/D/E/F/G/H/K/L/O)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B
$ws.Cells.Item(2, 2).Value = 0.9295119618572301
$ws.Cells.Item(3, 2).Value = 0.9018167476027656
$ws.Cells.Item(4, 2).Value = 0.8852098640543602
$ws.Cells.Item(5, 2).Value = 0.8785430886507584
$ws.Cells.Item(6, 2).Value = 0.8774421721109036
$ws.Cells.Item(7, 2).Value = 0.8851195452563729
$ws.Cells.Item(8, 2).Value = 0.9198803544065868
$ws.Cells.Item(9, 2).Value = 0.9911845110700312
$ws.Cells.Item(10, 2).Value = 1.045463893330748
$ws.Cells.Item(11, 2).Value = 1.070563818952252
$ws.Cells.Item(12, 2).Value = 1.080126694055622
$ws.Cells.Item(13, 2).Value = 1.078064584097291
$ws.Cells.Item(14, 2).Value = 1.071349401829849
$ws.Cells.Item(15, 2).Value = 1.06724370486927
$ws.Cells.Item(16, 2).Value = 1.043831713259806
$ws.Cells.Item(17, 2).Value = 1.029573286188366
$ws.Cells.Item(18, 2).Value = 1.021410650428265
$ws.Cells.Item(19, 2).Value = 1.01865354123413
$ws.Cells.Item(20, 2).Value = 1.03108714652879
$ws.Cells.Item(21, 2).Value = 1.073320242090574
$ws.Cells.Item(22, 2).Value = 1.101260262780272
$ws.Cells.Item(23, 2).Value = 1.08631739151798
$ws.Cells.Item(24, 2).Value = 1.030402622063633
$ws.Cells.Item(25, 2).Value = 0.9715610417142955

# Column D
$ws.Cells.Item(2, 4).Value = 0.02829632356917955
$ws.Cells.Item(3, 4).Value = 0.02555375907647317
$ws.Cells.Item(4, 4).Value = 0.02385840434255471
$ws.Cells.Item(5, 4).Value = 0.02316470416609917
$ws.Cells.Item(6, 4).Value = 0.02304934605886899
$ws.Cells.Item(7, 4).Value = 0.02384906026289713
$ws.Cells.Item(8, 4).Value = 0.02735307821608046
$ws.Cells.Item(9, 4).Value = 0.03413248711799355
$ws.Cells.Item(10, 4).Value = 0.03905577690838413
$ws.Cells.Item(11, 4).Value = 0.04128272609065675
$ws.Cells.Item(12, 4).Value = 0.04212415569094219
$ws.Cells.Item(13, 4).Value = 0.04194302254502702
$ws.Cells.Item(14, 4).Value = 0.04135198869605006
$ws.Cells.Item(15, 4).Value = 0.04098971884903335
$ws.Cells.Item(16, 4).Value = 0.03890998159604919
$ws.Cells.Item(17, 4).Value = 0.03763085093656571
$ws.Cells.Item(18, 4).Value = 0.03689393795845319
$ws.Cells.Item(19, 4).Value = 0.03664422870671302
$ws.Cells.Item(20, 4).Value = 0.03776714014870919
$ws.Cells.Item(21, 4).Value = 0.04152564070766118
$ws.Cells.Item(22, 4).Value = 0.04397113099017247
$ws.Cells.Item(23, 4).Value = 0.04266693945625377
$ws.Cells.Item(24, 4).Value = 0.03770552853884368
$ws.Cells.Item(25, 4).Value = 0.03230847543098037

# Column E
$ws.Cells.Item(2, 5).Value = 0.4073428573611526
$ws.Cells.Item(3, 5).Value = 0.4121655656644716
$ws.Cells.Item(4, 5).Value = 0.4153281861993481
$ws.Cells.Item(5, 5).Value = 0.4166676690771425
$ws.Cells.Item(6, 5).Value = 0.416893151611589
$ws.Cells.Item(7, 5).Value = 0.4153460456529832
$ws.Cells.Item(8, 5).Value = 0.4089639335475823
$ws.Cells.Item(9, 5).Value = 0.3980456364031486
$ws.Cells.Item(10, 5).Value = 0.3909953958569723
$ws.Cells.Item(11, 5).Value = 0.3879985390480201
$ws.Cells.Item(12, 5).Value = 0.3868939226029546
$ws.Cells.Item(13, 5).Value = 0.3871304773886735
$ws.Cells.Item(14, 5).Value = 0.3879070558699542
$ws.Cells.Item(15, 5).Value = 0.3883866686994182
$ws.Cells.Item(16, 5).Value = 0.3911954784419311
$ws.Cells.Item(17, 5).Value = 0.3929724511261643
$ws.Cells.Item(18, 5).Value = 0.3940143173880362
$ws.Cells.Item(19, 5).Value = 0.3943704763594891
$ws.Cells.Item(20, 5).Value = 0.3927812403991737
$ws.Cells.Item(21, 5).Value = 0.3876781356997121
$ws.Cells.Item(22, 5).Value = 0.3845191467659088
$ws.Cells.Item(23, 5).Value = 0.3861890434110045
$ws.Cells.Item(24, 5).Value = 0.3928676236566968
$ws.Cells.Item(25, 5).Value = 0.4008285768318114

# Column F
$ws.Cells.Item(2, 6).Value = 0.5540613752812575
$ws.Cells.Item(3, 6).Value = 0.5540019758901806
$ws.Cells.Item(4, 6).Value = 0.5543424821304797
$ws.Cells.Item(5, 6).Value = 0.5545760369707438
$ws.Cells.Item(6, 6).Value = 0.5546205440457754
$ws.Cells.Item(7, 6).Value = 0.5543452481065856
$ws.Cells.Item(8, 6).Value = 0.5539626520722223
$ws.Cells.Item(9, 6).Value = 0.5562043268127113
$ws.Cells.Item(10, 6).Value = 0.5596776131576746
$ws.Cells.Item(11, 6).Value = 0.5616548401247812
$ws.Cells.Item(12, 6).Value = 0.5624607022612835
$ws.Cells.Item(13, 6).Value = 0.5622846044436187
$ws.Cells.Item(14, 6).Value = 0.5617199938785902
$ws.Cells.Item(15, 6).Value = 0.5613815937996307
$ws.Cells.Item(16, 6).Value = 0.5595563905871543
$ws.Cells.Item(17, 6).Value = 0.5585384310647683
$ws.Cells.Item(18, 6).Value = 0.5579903145754486
$ws.Cells.Item(19, 6).Value = 0.5578111526353808
$ws.Cells.Item(20, 6).Value = 0.5586429251069447
$ws.Cells.Item(21, 6).Value = 0.5618842831414597
$ws.Cells.Item(22, 6).Value = 0.5643356944134297
$ws.Cells.Item(23, 6).Value = 0.5629968578044924
$ws.Cells.Item(24, 6).Value = 0.5585955677195784
$ws.Cells.Item(25, 6).Value = 0.5552773367387118

# Column G
$ws.Cells.Item(2, 7).Value = 0.394618540327393
$ws.Cells.Item(3, 7).Value = 0.3959734702752016
$ws.Cells.Item(4, 7).Value = 0.3971155826816357
$ws.Cells.Item(5, 7).Value = 0.3976589189305386
$ws.Cells.Item(6, 7).Value = 0.3977538430560585
$ws.Cells.Item(7, 7).Value = 0.3971225949332151
$ws.Cells.Item(8, 7).Value = 0.3950212928324461
$ws.Cells.Item(9, 7).Value = 0.3933659107447838
$ws.Cells.Item(10, 7).Value = 0.3936589791832006
$ws.Cells.Item(11, 7).Value = 0.394121428441693
$ws.Cells.Item(12, 7).Value = 0.3943439729855669
$ws.Cells.Item(13, 7).Value = 0.3942939332753213
$ws.Cells.Item(14, 7).Value = 0.3941387862997914
$ws.Cells.Item(15, 7).Value = 0.3940499332742888
$ws.Cells.Item(16, 7).Value = 0.3936353882688479
$ws.Cells.Item(17, 7).Value = 0.3934654457461306
$ws.Cells.Item(18, 7).Value = 0.3933986729955166
$ws.Cells.Item(19, 7).Value = 0.3933813816894727
$ws.Cells.Item(20, 7).Value = 0.3934803300534213
$ws.Cells.Item(21, 7).Value = 0.3941830688641375
$ws.Cells.Item(22, 7).Value = 0.394918822442861
$ws.Cells.Item(23, 7).Value = 0.3945008081326193
$ws.Cells.Item(24, 7).Value = 0.3934735045225608
$ws.Cells.Item(25, 7).Value = 0.3935491467586374

# Column H
$ws.Cells.Item(2, 8).Value = 0.5598256522037914
$ws.Cells.Item(3, 8).Value = 0.5639930804856448
$ws.Cells.Item(4, 8).Value = 0.5668150632768132
$ws.Cells.Item(5, 8).Value = 0.5680312582803992
$ws.Cells.Item(6, 8).Value = 0.5682372068527854
$ws.Cells.Item(7, 8).Value = 0.5668311971578817
$ws.Cells.Item(8, 8).Value = 0.5612079849515368
$ws.Cells.Item(9, 8).Value = 0.5522673742740665
$ws.Cells.Item(10, 8).Value = 0.5469687022828396
$ws.Cells.Item(11, 8).Value = 0.5448335629391465
$ws.Cells.Item(12, 8).Value = 0.5440645891332139
$ws.Cells.Item(13, 8).Value = 0.5442284424698869
$ws.Cells.Item(14, 8).Value = 0.5447695063548537
$ws.Cells.Item(15, 8).Value = 0.5451060744003371
$ws.Cells.Item(16, 8).Value = 0.5471137754131377
$ws.Cells.Item(17, 8).Value = 0.548415916885034
$ws.Cells.Item(18, 8).Value = 0.5491907837820094
$ws.Cells.Item(19, 8).Value = 0.549457591121012
$ws.Cells.Item(20, 8).Value = 0.5482746203324069
$ws.Cells.Item(21, 8).Value = 0.5446095093987964
$ws.Cells.Item(22, 8).Value = 0.5424447055428345
$ws.Cells.Item(23, 8).Value = 0.54357901396736
$ws.Cells.Item(24, 8).Value = 0.5483384187396894
$ws.Cells.Item(25, 8).Value = 0.554462878284987

# Column K
$ws.Cells.Item(2, 11).Value = 0.5065805752580843
$ws.Cells.Item(3, 11).Value = 0.4556783212029529
$ws.Cells.Item(4, 11).Value = 0.4242478554418767
$ws.Cells.Item(5, 11).Value = 0.411396273404705
$ws.Cells.Item(6, 11).Value = 0.4092596830052742
$ws.Cells.Item(7, 11).Value = 0.4240747090709931
$ws.Cells.Item(8, 11).Value = 0.4890665960867011
$ws.Cells.Item(9, 11).Value = 0.6150825955405139
$ws.Cells.Item(10, 11).Value = 0.7067557507934055
$ws.Cells.Item(11, 11).Value = 0.7482551600473357
$ws.Cells.Item(12, 11).Value = 0.7639398896737646
$ws.Cells.Item(13, 11).Value = 0.7605632610429325
$ws.Cells.Item(14, 11).Value = 0.749546162963469
$ws.Cells.Item(15, 11).Value = 0.7427939108582393
$ws.Cells.Item(16, 11).Value = 0.7040394945763353
$ws.Cells.Item(17, 11).Value = 0.6802121795777794
$ws.Cells.Item(18, 11).Value = 0.666488273245335
$ws.Cells.Item(19, 11).Value = 0.6618383499150298
$ws.Cells.Item(20, 11).Value = 0.682750618207649
$ws.Cells.Item(21, 11).Value = 0.7527829790260796
$ws.Cells.Item(22, 11).Value = 0.7983767186534863
$ws.Cells.Item(23, 11).Value = 0.7740589501051431
$ws.Cells.Item(24, 11).Value = 0.6816030687286627
$ws.Cells.Item(25, 11).Value = 0.5811492043175406

# Column L
$ws.Cells.Item(2, 12).Value = 0.1227057796620556
$ws.Cells.Item(3, 12).Value = 0.1113939134944673
$ws.Cells.Item(4, 12).Value = 0.1044729462433054
$ws.Cells.Item(5, 12).Value = 0.1016589075174466
$ws.Cells.Item(6, 12).Value = 0.1011920237133523
$ws.Cells.Item(7, 12).Value = 0.1044349693981843
$ws.Cells.Item(8, 12).Value = 0.1188004270693597
$ws.Cells.Item(9, 12).Value = 0.1471614772530074
$ws.Cells.Item(10, 12).Value = 0.168110547274722
$ws.Cells.Item(11, 12).Value = 0.1776644925446647
$ws.Cells.Item(12, 12).Value = 0.181285682574952
$ws.Cells.Item(13, 12).Value = 0.1805056492910069
$ws.Cells.Item(14, 12).Value = 0.1779623443741656
$ws.Cells.Item(15, 12).Value = 0.1764049244132195
$ws.Cells.Item(16, 12).Value = 0.1674866469123231
$ws.Cells.Item(17, 12).Value = 0.1620216415938103
$ws.Cells.Item(18, 12).Value = 0.1588805931303767
$ws.Cells.Item(19, 12).Value = 0.1578174838722504
$ws.Cells.Item(20, 12).Value = 0.1626031660458409
$ws.Cells.Item(21, 12).Value = 0.1787092858034356
$ws.Cells.Item(22, 12).Value = 0.1892548389055975
$ws.Cells.Item(23, 12).Value = 0.183624764314132
$ws.Cells.Item(24, 12).Value = 0.1623402561649101
$ws.Cells.Item(25, 12).Value = 0.139469077012663

# Column O
$ws.Cells.Item(2, 15).Value = 1.853650692566276
$ws.Cells.Item(3, 15).Value = 1.865005279891861
$ws.Cells.Item(4, 15).Value = 1.873177042429106
$ws.Cells.Item(5, 15).Value = 1.876808755959516
$ws.Cells.Item(6, 15).Value = 1.877430015551383
$ws.Cells.Item(7, 15).Value = 1.873224799731332
$ws.Cells.Item(8, 15).Value = 1.857316634766292
$ws.Cells.Item(9, 15).Value = 1.835647885936794
$ws.Cells.Item(10, 15).Value = 1.825545589339669
$ws.Cells.Item(11, 15).Value = 1.822215259086931
$ws.Cells.Item(12, 15).Value = 1.8211362369764
$ws.Cells.Item(13, 15).Value = 1.821360522170096
$ws.Cells.Item(14, 15).Value = 1.822122836798087
$ws.Cells.Item(15, 15).Value = 1.822613496208703
$ws.Cells.Item(16, 15).Value = 1.825788708137935
$ws.Cells.Item(17, 15).Value = 1.828060774609924
$ws.Cells.Item(18, 15).Value = 1.829486689630116
$ws.Cells.Item(19, 15).Value = 1.829989927571887
$ws.Cells.Item(20, 15).Value = 1.827806583498131
$ws.Cells.Item(21, 15).Value = 1.821893983071561
$ws.Cells.Item(22, 15).Value = 1.81909125642315
$ws.Cells.Item(23, 15).Value = 1.820489945291911
$ws.Cells.Item(24, 15).Value = 1.827921130536566
$ws.Cells.Item(25, 15).Value = 1.840488894684498
